$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly block (rows 291-292),
# pushing the existing data down by two rows.
$ws.Rows.Item(291).Insert()
$ws.Rows.Item(292).Insert()

# Row 291: Apio, Americana (o), Primera - new week (2022-01-07 / serial 44568)
$ws.Cells.Item(291, 1).Value = 8
$ws.Cells.Item(291, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(291, 3).Value = "Coquimbo"
$ws.Cells.Item(291, 4).Value = 44568
$ws.Cells.Item(291, 5).Value = 4
$ws.Cells.Item(291, 6).Value = 100112017
$ws.Cells.Item(291, 7).Value = "Apio"
$ws.Cells.Item(291, 8).Value = "Americana (o)"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 2560
$ws.Cells.Item(291, 11).Value = 8000
$ws.Cells.Item(291, 12).Value = 9000
$ws.Cells.Item(291, 13).Value = 8500
$ws.Cells.Item(291, 14).Value = "`$/docena de matas"
$ws.Cells.Item(291, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(291, 16).Value = 1417
$ws.Cells.Item(291, 17).Value = 6
$ws.Cells.Item(291, 18).Value = "Hortaliza"

# Row 292: Apio, Americana (o), Segunda - new week (2022-01-07 / serial 44568)
$ws.Cells.Item(292, 1).Value = 8
$ws.Cells.Item(292, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(292, 3).Value = "Coquimbo"
$ws.Cells.Item(292, 4).Value = 44568
$ws.Cells.Item(292, 5).Value = 4
$ws.Cells.Item(292, 6).Value = 100112017
$ws.Cells.Item(292, 7).Value = "Apio"
$ws.Cells.Item(292, 8).Value = "Americana (o)"
$ws.Cells.Item(292, 9).Value = "Segunda"
$ws.Cells.Item(292, 10).Value = 1500
$ws.Cells.Item(292, 11).Value = 6000
$ws.Cells.Item(292, 12).Value = 7000
$ws.Cells.Item(292, 13).Value = 6500
$ws.Cells.Item(292, 14).Value = "`$/docena de matas"
$ws.Cells.Item(292, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(292, 16).Value = 1083
$ws.Cells.Item(292, 17).Value = 6
$ws.Cells.Item(292, 18).Value = "Hortaliza"

Write-Output "done"
